$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.018.93'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '3.568.70'
$ws.Range("E3").Value = '  +2.48%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '601.13'
$ws.Range("E5").Value = '  +2.03%  '
$ws.Range("D6").Value = '135.63'
$ws.Range("E6").Value = '  -1.83%  '
$ws.Range("D7").Value = '3.567.43'
$ws.Range("E7").Value = '  +2.36%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("E10").Value = '  +0.42%  '
$ws.Range("D11").Value = '6.94'
$ws.Range("E11").Value = '  -4.12%  '
$ws.Range("E12").Value = '  +1.37%  '
$ws.Range("D13").Value = '4.182.28'
$ws.Range("E13").Value = '  +2.75%  '
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("D15").Value = '3.574.74'
$ws.Range("E15").Value = '  +2.78%  '
$ws.Range("D16").Value = '27.14'
$ws.Range("E16").Value = '  +1.86%  '
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("D18").Value = '65.146.90'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = '10.09'
$ws.Range("E19").Value = '  +3.72%  '
$ws.Range("E20").Value = '  +3.91%  '
$ws.Range("E21").Value = '  +1.23%  '
$ws.Range("D22").Value = '388.46'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  +4.41%  '
$ws.Range("D24").Value = '3.719.63'
$ws.Range("E24").Value = '  +2.70%  '
$ws.Range("D25").Value = '74.19'
$ws.Range("E25").Value = '  +2.21%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '0.0000116'
$ws.Range("E27").Value = '  +5.33%  '
$ws.Range("E28").Value = '  +3.98%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").Value = '2.29'
$ws.Range("E30").Value = '  +3.29%  '
$ws.Range("D31").Value = '8.41'
$ws.Range("E31").Value = '  +2.05%  '
$ws.Range("E32").Value = '  +21.58%  '
$ws.Range("D33").Value = '3.583.81'
$ws.Range("E33").Value = '  +2.50%  '
$ws.Range("D34").Value = '24.02'
$ws.Range("E34").Value = '  +4.20%  '
$ws.Range("E36").Value = '  +1.28%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '169.67'
$ws.Range("E37").Value = '  -1.39%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").Value = '6.93'
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("E39").Value = '  +5.66%  '
$ws.Range("D40").Value = '5.01'
$ws.Range("E40").Value = '  +5.40%  '
$ws.Range("E41").Value = '  +3.92%  '
$ws.Range("D42").Value = '27.25'
$ws.Range("E42").Value = '  +10.52%  '
$ws.Range("D43").Value = '0.825'
$ws.Range("E43").Value = '  +1.46%  '
$ws.Range("D44").Value = '42.67'
$ws.Range("E44").Value = '  +0.67%  '
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("E46").Value = '  +2.82%  '
$ws.Range("E47").Value = '  +5.17%  '
$ws.Range("D48").Value = '1.65'
$ws.Range("E48").Value = '  +2.13%  '
$ws.Range("D49").Value = '2.505.10'
$ws.Range("E49").Value = '  +12.24%  '
$ws.Range("D50").Value = '6.93'
$ws.Range("E50").Value = '  +3.89%  '
$ws.Range("E51").Value = '  +10.21%  '
